# ODBC2KML Requirements Meeting - bump the meeting date from
# October 27, 2009 to November 5, 2009 (commit: "Changed the date to
# November 5, 2009"), and tidy up a previously-split bullet on the
# "Questions for the Client" slide.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Title slide ("Subtitle 2"): the big visible "October 27, 2009"
#    becomes "November 5, 2009". Only the "October 27, " portion is
#    retyped -- "2009" is left as-is, which is why the saved file ends
#    up with two runs in that paragraph instead of one.
# ---------------------------------------------------------------------
$titleSlide = $p.Slides.Item(1)
$subtitle = $titleSlide.Shapes.Item(2)
$subtitleText = $subtitle.TextFrame.TextRange
$fullText = $subtitleText.Text
$oldDate = "October 27, 2009"
$pos = $fullText.IndexOf($oldDate) + 1
if ($pos -gt 0) {
    $datePrefix = $subtitleText.Characters($pos, 12)   # "October 27, "
    $datePrefix.Text = "November 5, "
}

# ---------------------------------------------------------------------
# 2. Auto "Last saved"/"Today" date placeholders. These live on the
#    slide master, every slide layout, and the notes master, each
#    caching the field as 10/29/2009 -> refresh them all to 11/5/2009.
# ---------------------------------------------------------------------
function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "10/29/2009") {
                $tr.Text = "11/5/2009"
            }
        }
    }
}

Update-DatePlaceholder $p.SlideMaster.Shapes

$customLayouts = $p.SlideMaster.CustomLayouts
for ($L = 1; $L -le $customLayouts.Count; $L++) {
    Update-DatePlaceholder $customLayouts.Item($L).Shapes
}

# The notes master's "Date Placeholder" shape shares its internal shape
# id with the slide master's footer placeholder, so reaching it through
# Shapes.Item(...) cross-talks with the slide master. Go through the
# notes master's own HeadersFooters/DateAndTime instead, which resolves
# correctly and updates the notes master in place.
$p.NotesMaster.HeadersFooters.DateAndTime.Text = "11/5/2009"

# ---------------------------------------------------------------------
# 3. "Questions for the Client" slide: the bullet "Individual
#    descriptions for table rows" used to be split across two runs
#    ("Individual descriptions for table " + "rows"); merge it back
#    into a single run.
# ---------------------------------------------------------------------
$questionsSlide = $p.Slides.Item(11)
$content = $questionsSlide.Shapes.Item(2)
$contentText = $content.TextFrame.TextRange
for ($i = 1; $i -le $contentText.Paragraphs().Count; $i++) {
    $para = $contentText.Paragraphs($i)
    if ($para.Text -eq "Individual descriptions for table rows") {
        $para.Text = "TEMP_MERGE_PLACEHOLDER"
        $para2 = $contentText.Paragraphs($i)
        $para2.Text = "Individual descriptions for table rows"
        break
    }
}
